$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.834.41"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").Value = "3.785.99"
$ws.Range("E3").Value = "  +20.50%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'615.50"
$ws.Range("E5").Value = "  +6.26%  "

$ws.Range("D6").Value = "'177.44"
$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").Value = "3.776.57"
$ws.Range("E7").Value = "  +20.30%  "

$ws.Range("E9").Value = "  +4.38%  "

$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +9.60%  "

$ws.Range("E11").Value = "  -2.78%  "

$ws.Range("D12").Value = "'0.501"
$ws.Range("E12").Value = "  +5.87%  "

$ws.Range("D13").Value = "'40.66"
$ws.Range("E13").Value = "  +9.87%  "

$ws.Range("D14").Value = "'0.0000257"
$ws.Range("E14").Value = "  +5.40%  "

$ws.Range("D15").Value = "4.419.99"
$ws.Range("E15").Value = "  +20.64%  "

$ws.Range("D16").Value = "3.787.67"
$ws.Range("E16").Value = "  +20.88%  "

$ws.Range("D17").Value = "70.029.96"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  +6.58%  "

$ws.Range("D20").Value = "'518.20"
$ws.Range("E20").Value = "  +5.74%  "

$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").Value = "'9.44"
$ws.Range("E22").Value = "  +21.29%  "

$ws.Range("E23").Value = "  +5.85%  "

$ws.Range("D24").Value = "'88.97"
$ws.Range("E24").Value = "  +5.85%  "

$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +5.63%  "

$ws.Range("D26").Value = "'13.59"
$ws.Range("E26").Value = "  +5.22%  "

$ws.Range("D27").Value = "'10.89"
$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("D28").Value = "'0.0000128"
$ws.Range("E28").Value = "  +33.77%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +6.08%  "

$ws.Range("E31").Value = "  +8.59%  "

$ws.Range("E32").Value = "  -3.49%  "

$ws.Range("D33").Value = "'32.16"
$ws.Range("E33").Value = "  +13.24%  "

$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").Value = "'6.23"
$ws.Range("E36").Value = "  +9.55%  "

$ws.Range("E37").Value = "  +9.20%  "

$ws.Range("E38").Value = "  +5.26%  "

$ws.Range("E39").Value = "  +5.79%  "

$ws.Range("E40").Value = "  +5.97%  "

$ws.Range("D41").Value = "'51.52"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").Value = "'44.65"
$ws.Range("E42").Value = "  -8.99%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.84"
$ws.Range("E43").Value = "  +4.98%  "

$ws.Range("D44").Value = "'426.65"
$ws.Range("E44").Value = "  +9.25%  "

$ws.Range("D45").Value = "3.115.49"
$ws.Range("E45").Value = "  +11.70%  "

$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("E47").Value = "  +4.45%  "

$ws.Range("D48").Value = "'27.80"
$ws.Range("E48").Value = "  +2.89%  "

$ws.Range("D49").Value = "'137.44"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("E51").Value = "  +5.82%  "
